# Apply the edits described by the diff to RS-BPV.schema.docx
#
# Find.Execute with a non-empty ReplaceWith performs a document-wide
# replace (it is not confined to the calling Range), and it also runs
# the replacement text through Word's AutoCorrect/"smart quotes"
# machinery. To keep each edit scoped to exactly the paragraph/cell it
# belongs to, and to avoid mangling straight apostrophes into curly
# ones, we instead use Find.Execute purely to *locate* the (unique)
# target text - passing an empty ReplaceWith and Replace:=wdReplaceNone
# (0) - and then assign the resulting Range's .Text directly.

$d = $word.ActiveDocument

function Replace-UniqueText($range, [string]$oldText, [string]$newText) {
    $found = $range.Find.Execute(
        $oldText, $true, $false, $false, $false, $false,
        $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $oldText"
    }
    $range.Text = $newText
}

# 1) caseId description: shortened text (drop the "de régulation médicale (DRM)" part)
Replace-UniqueText $d.Content `
    "Identifiant partagé du dossier de régulation médicale (DRM)" `
    "Identifiant partagé du dossier "

# 2) "A valoriser avec le prénom..." description: extend the sentence
Replace-UniqueText $d.Content `
    "A valoriser avec le prénom et le nom du rédacteur ou un numéro RPPS. " `
    "A valoriser avec le prénom et le nom du rédacteur, un numéro RPPS, un matricule, etc. "

# 3) "procedure" field description: "le SMUR" -> "la ressource"
Replace-UniqueText $d.Content `
    "Actes réalisés par le SMUR" `
    "Actes réalisés par la ressource"

# 4) "associatedDiagnosis" cardinality: 0..1 -> 0..n (Table 4, row 4, column 4)
$cell = $d.Tables.Item(4).Cell(4, 4)
$cell.Range.Text = "0..n"

# 5) "freetext" field description: fix double comma typo
Replace-UniqueText $d.Content `
    "Permettrait de concaténer dans une zone de commentaire d'autres champs (ex. anamnèse : allergies,, traitements, symptomes, antécédents)" `
    "Permettrait de concaténer dans une zone de commentaire d'autres champs (ex. anamnèse : allergies, traitements, symptomes, antécédents)"
